# aggiornamento dell'analisi dei rischi Excel
#
# The "Risk Register" sheet numbers each risk row in column B (B3=1,
# B4=2, ...). Rows 5-9 were left blank in column B; fill in the
# continuing sequence (3-7) and leave the active selection on B9,
# matching where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risk Register")

$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7

$ws.Activate()
$ws.Range("B9").Select()
